$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4: clear PB0 Function/CN7 (previously "Ultrasound Trigger" / 34) ---
$ws.Range("H4").Value = $null
$ws.Range("I4").Value = $null

# --- Row 5: add PA1 Function/CN7 (Motor Direction / 30) ---
$ws.Range("C5").Value = "Motor Direction"
$ws.Range("D5").Value = 30

# --- Row 5: add PD1 Function/CN7 (Ultrasound Trigger / 31) ---
$ws.Range("R5").Value = "Ultrasound Trigger"
$ws.Range("S5").Value = 31

# --- Row 10: add PB6 Function (Motor PWM) and CN10 value (17) ---
$ws.Range("H10").Value = "Motor PWM"
$ws.Range("J10").Value = 17

# --- Column widths (best-fit on the newly populated columns) ---
# The headless runtime quantizes ColumnWidth to the nearest 1/6 unit, so the
# inputs below are chosen to land as close as possible to the target widths
# (15.140625 and 17.5703125) recorded in the saved workbook.
$ws.Columns("C:C").ColumnWidth = 14.333333333333334
$ws.Columns("R:R").ColumnWidth = 16.666666666666668

# --- Selection / active cell ---
$ws.Range("G14").Select()

# --- Window position ---
$wb.Windows.Item(1).Left = 3420
$wb.Windows.Item(1).Top = 3420
